# Update "想去人数" (number of people interested) values in column F
# for the "展览" (Exhibition) and "全部类型" (All types) worksheets,
# reflecting a refreshed data scrape (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14044
$ws1.Range("F4").Value = 547
$ws1.Range("F6").Value = 1216
$ws1.Range("F7").Value = 1044
$ws1.Range("F8").Value = 13920
$ws1.Range("F9").Value = 14928
$ws1.Range("F11").Value = 11
$ws1.Range("F21").Value = 21
$ws1.Range("F22").Value = 1153
$ws1.Range("F25").Value = 5810
$ws1.Range("F28").Value = 5440
$ws1.Range("F29").Value = 56
$ws1.Range("F30").Value = 128
$ws1.Range("F31").Value = 63
$ws1.Range("F32").Value = 316

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14044
$ws4.Range("F5").Value = 547
$ws4.Range("F7").Value = 1216
$ws4.Range("F8").Value = 1044
$ws4.Range("F9").Value = 13920
$ws4.Range("F10").Value = 14928
$ws4.Range("F12").Value = 11
$ws4.Range("F22").Value = 21
$ws4.Range("F23").Value = 1153
$ws4.Range("F27").Value = 5810
$ws4.Range("F30").Value = 5440
$ws4.Range("F31").Value = 56
$ws4.Range("F32").Value = 128
$ws4.Range("F33").Value = 63
$ws4.Range("F34").Value = 316
